# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Swap "Santa Lucia" / "Timor Oriental" labels (rows 206 & 207) ---
# Before: A206 = "Santa Lucia", A207 = "Timor Oriental"
# After:  A206 = "Timor Oriental", A207 = "Santa Lucia"
$ws.Range("A206").Value = "Timor Oriental"
$ws.Range("A207").Value = "Santa Lucia"

# --- Update "last updated" timestamp text (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 27 de Septiembre de 2020 a las 05:58"

# --- Update country statistics ---
# Row 5: India
$ws.Range("B5").Value = 5992532
$ws.Range("C5").Value = 1951
$ws.Range("D5").Value = 4941627
$ws.Range("E5").Value = 956371

# Row 35: Belgica
$ws.Range("B35").Value = 112803
$ws.Range("C35").Value = 1827
$ws.Range("D35").Value = 19246
$ws.Range("E35").Value = 83583
$ws.Range("G35").Value = 5
$ws.Range("H35").Value = 9974

# Row 39: Kazajistan
$ws.Range("D39").Value = 102666
$ws.Range("E39").Value = 3358

# Row 50: Honduras
$ws.Range("B50").Value = 74548
$ws.Range("C50").Value = 708
$ws.Range("D50").Value = 26088
$ws.Range("E50").Value = 46172
$ws.Range("G50").Value = 17
$ws.Range("H50").Value = 2288

# Row 78: Australia
$ws.Range("B78").Value = 27040
$ws.Range("C78").Value = 24
$ws.Range("D78").Value = 24573
$ws.Range("E78").Value = 1595

# Row 158: Belice
$ws.Range("B158").Value = 1825
$ws.Range("C158").Value = 17
$ws.Range("D158").Value = 1165
$ws.Range("E158").Value = 636
$ws.Range("G158").Value = 1
$ws.Range("H158").Value = 24

# Row 187: Camboya
$ws.Range("B187").Value = 276
$ws.Range("C187").Value = 1
$ws.Range("E187").Value = 2

# Row 188: Butan
$ws.Range("B188").Value = 271
$ws.Range("C188").Value = 8
$ws.Range("D188").Value = 205
$ws.Range("E188").Value = 66
